# update made january 24th
# Re-sorts the CVR/Lost TCV table by CVR ascending and adds 11 new CVR rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("10082366", 4512, 1),
    @("10750113", 23352, 1),
    @("12530633", 128134, 1),
    @("13893179", 40172, 1),
    @("14364641", 45106, 1),
    @("16049441", 199392.4, 1),
    @("17038117", 125865, 1),
    @("18788799", 116330, 1),
    @("19029484", 80135, 1),
    @("19195597", 207498.12, 1),
    @("19248615", 922780, 1),
    @("20039701", 51093.54, 1),
    @("20309687", 76044, 1),
    @("21699780", 21650.08, 1),
    @("25086589", 278110, 1),
    @("26317797", 132000, 1),
    @("26688205", 74972.74000000001, 1),
    @("27276407", 75401.86, 1),
    @("27351034", 71414.8, 1),
    @("27524958", 163235.74, 1),
    @("28979355", 20119, 1),
    @("29186898", 23777, 1),
    @("30715993", 6371, 1),
    @("32440975", 6137, 1),
    @("32775705", 74171, 1),
    @("33372477", 70780, 1),
    @("33738811", 105096, 1),
    @("34082200", 172800, 1),
    @("34725101", 30453, 1),
    @("35432604", 128519, 1),
    @("35802495", 171628, 1),
    @("36184779", 140432, 1),
    @("36934085", 129827.5, 1),
    @("38618873", 98144, 1),
    @("42268992", 32723, 1),
    @("46936728", 160611, 1),
    @("47252059", 31857, 1),
    @("55164711", 46106, 1),
    @("67437853", 30097, 1),
    @("75128711", 4513, 1),
    @("87719715", 258702.76, 1)
)

$rowCount = $data.Length
$lastRow = $rowCount + 1

# Clear out any old rows below the new data range (table shrank never happens here,
# but keep this safe in case old content lingers beyond the new extent).
$ws.Range("A2:C10000").ClearContents()

# Column A (CVR) must stay text, matching the original shared-string storage,
# so force text formatting before writing the values, then strip the
# formatting back off so no extra style index lingers on the cells.
$ws.Range("A2:A$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

$ws.Range("A2:A$lastRow").ClearFormats()
